$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 385 (shifts existing data rows 385-415 down to 386-416)
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new weekly price record
$ws.Cells.Item(385,1).Value  = 4
$ws.Cells.Item(385,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(385,3).Value  = "Los Lagos"
$ws.Cells.Item(385,4).Value  = 45106
$ws.Cells.Item(385,5).Value  = 10
$ws.Cells.Item(385,6).Value  = 100112032
$ws.Cells.Item(385,7).Value  = "Zapallo italiano"
$ws.Cells.Item(385,8).Value  = "Sin especificar"
$ws.Cells.Item(385,9).Value  = "Primera"
$ws.Cells.Item(385,10).Value = 140
$ws.Cells.Item(385,11).Value = 18000
$ws.Cells.Item(385,12).Value = 19000
$ws.Cells.Item(385,13).Value = 18500
$ws.Cells.Item(385,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(385,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(385,16).Value = 370
$ws.Cells.Item(385,17).Value = 50
$ws.Cells.Item(385,18).Value = "Hortaliza"
